# default types sheet updated
# Insert two new "default type" rows — "Site Group" and "People Group" —
# right after the existing "Pgroup Identifier" row (row 5), pushing every
# row below it down by two. Both new rows reference "PGROUPIDENTIFIER" as
# their Type (column D), matching the other rows that hang off Pgroup
# Identifier-like parents, with BV/Client left as "NONE".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 6 (formatting/style is inherited from
# the row above, same as Excel's native "Insert Copied Cells" behaviour).
$ws.Rows("6:7").Insert()
$ws.Rows("6:7").RowHeight = 18.75

# Row 6: Site Group
$ws.Range("B6").Value = "Site Group"
$ws.Range("C6").Value = "SITEGROUP"
$ws.Range("D6").Value = "PGROUPIDENTIFIER"
$ws.Range("E6").Value = "NONE"
$ws.Range("F6").Value = "NONE"

# Row 7: People Group
$ws.Range("B7").Value = "People Group"
$ws.Range("C7").Value = "PEOPLEGROUP"
$ws.Range("D7").Value = "PGROUPIDENTIFIER"
$ws.Range("E7").Value = "NONE"
$ws.Range("F7").Value = "NONE"

# Restore the view to match the saved workbook state (active cell J6
# instead of the pre-edit D20 selection further down the sheet).
$ws.Range("J6").Select()
